$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: I = "I0", J = "IF" (added after existing A:H columns)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match header formatting used by the other header cells (bold, centered,
# top-aligned, thin box border)
foreach ($cellRef in @("I1", "J1")) {
    $headerCell = $ws.Range($cellRef)
    $headerCell.Font.Bold = $true
    $headerCell.HorizontalAlignment = -4108
    $headerCell.VerticalAlignment = -4160
    $headerCell.Borders.LineStyle = 1
}

# Per-row values for the new I0 / IF columns
$newColumnData = @(
    @{Row=2; I0=4; IF=5},
    @{Row=3; I0=8; IF=8},
    @{Row=4; I0=1; IF=3},
    @{Row=5; I0=8; IF=9},
    @{Row=6; I0=6; IF=7},
    @{Row=7; I0=1; IF=2},
    @{Row=8; I0=4; IF=6},
    @{Row=9; I0=7; IF=8},
    @{Row=10; I0=6; IF=6},
    @{Row=11; I0=9; IF=9},
    @{Row=12; I0=7; IF=7},
    @{Row=13; I0=1; IF=2},
    @{Row=14; I0=1; IF=2},
    @{Row=15; I0=7; IF=7},
    @{Row=16; I0=1; IF=2},
    @{Row=17; I0=7; IF=8},
    @{Row=18; I0=1; IF=2},
    @{Row=19; I0=10; IF=10},
    @{Row=20; I0=7; IF=8},
    @{Row=21; I0=7; IF=7},
    @{Row=22; I0=9; IF=9},
    @{Row=23; I0=8; IF=8},
    @{Row=24; I0=8; IF=8},
    @{Row=25; I0=9; IF=9},
    @{Row=26; I0=3; IF=3},
    @{Row=27; I0=1; IF=3},
    @{Row=28; I0=6; IF=6},
    @{Row=29; I0=1; IF=2},
    @{Row=30; I0=4; IF=5},
    @{Row=31; I0=9; IF=9},
    @{Row=32; I0=6; IF=7},
    @{Row=33; I0=7; IF=7},
    @{Row=34; I0=7; IF=8},
    @{Row=35; I0=7; IF=7},
    @{Row=36; I0=5; IF=5},
    @{Row=37; I0=5; IF=5}
)

foreach ($item in $newColumnData) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I0
    $ws.Cells.Item($item.Row, 10).Value = $item.IF
}
